$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18 (shifts existing rows 18-31 down to 19-32),
# copying formatting from the row above it so styles match the rest of the table.
$ws.Rows.Item(18).Insert(-4121, 0)   # xlShiftDown = -4121, xlFormatFromLeftOrAbove = 0

# Fill in the new item's data (item #12 in the list)
$ws.Cells.Item(18, 1).Value = 12                                        # A18 - sequence number
$ws.Cells.Item(18, 3).Value = "T4-THYRO 100MCG 100 TABLETS"             # C18 - item name
$ws.Cells.Item(18, 8).Value = "0:0"                                     # H18 - current balance
$ws.Cells.Item(18, 12).Value = "1"                                      # L18 - order limit
$ws.Cells.Item(18, 14).Value = "58.00"                                  # N18 - price
$ws.Cells.Item(18, 16).Value = "58.0000"                                # P18 - sale price
$ws.Cells.Item(18, 17).Value = "1:0"                                    # Q18 - transactions count

# Update the generated timestamp shown at the bottom of the report
$ws.Cells.Item(33, 1).Value = "Saturday, 13 September, 2025 12:38 PM"

# Update the total sum cell (P column) to include the new row's price
$ws.Cells.Item(32, 16).Value = 756.57000000000005
